$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.37898
$ws.Range("H2").Value = 1.13694
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.367566666666667
$ws.Range("N2").Value = 10.1027
$ws.Range("O2").Value = 0.9944910078726888
$ws.Range("P2").Value = 0.9944910078726888
$ws.Range("Q2").Value = 1.276240415333334
$ws.Range("R2").Value = 11.486163738
$ws.Range("S2").Value = 0.9944910078726888
$ws.Range("T2").Value = 0.9944910078726888

# --- Add row 3 (FAPs -> FAPs) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Nell2"
$ws.Range("C3").Value = "Robo3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.37898
$ws.Range("H3").Value = 1.13694
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01387666666666667
$ws.Range("N3").Value = 0.04163
$ws.Range("O3").Value = 0.004097979813093533
$ws.Range("P3").Value = 0.004097979813093532
$ws.Range("Q3").Value = 0.005258979133333334
$ws.Range("R3").Value = 0.0473308122
$ws.Range("S3").Value = 0.004097979813093533
$ws.Range("T3").Value = 0.004097979813093532

# --- Add row 4 (FAPs -> sCs) ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nell2"
$ws.Range("C4").Value = "Robo3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.37898
$ws.Range("H4").Value = 1.13694
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.004778
$ws.Range("N4").Value = 0.014334
$ws.Range("O4").Value = 0.001411012314217696
$ws.Range("P4").Value = 0.001411012314217696
$ws.Range("Q4").Value = 0.00181076644
$ws.Range("R4").Value = 0.01629689796
$ws.Range("S4").Value = 0.001411012314217696
$ws.Range("T4").Value = 0.001411012314217696
